# Quarterly indexing esoteric bug-fix operation
#
# Column A holds the "as-of" date serials for each forecast row. Due to an
# indexing bug, those dates were being stamped as the 1st of the quarter's
# start month. The fix re-stamps each row one month later, on the 15th (the
# mid-point of the following month), leaving every other column/value
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldSerial = $cell.Value2
    if ($oldSerial -eq $null) { continue }

    $oldDate = [DateTime]::FromOADate([double]$oldSerial)

    $newMonth = $oldDate.Month + 1
    $newYear = $oldDate.Year
    if ($newMonth -gt 12) {
        $newMonth = $newMonth - 12
        $newYear = $newYear + 1
    }

    $newDate = Get-Date -Year $newYear -Month $newMonth -Day 15 -Hour 0 -Minute 0 -Second 0
    $cell.Value2 = $newDate.ToOADate()
}
